$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example - Project Plan Template")

# ---------------------------------------------------------------------------
# Task status + date updates (rows 18-35 of the "Crono" schedule table)
# ---------------------------------------------------------------------------

# Row 18: 2.4 -> Complete, end date slips to 45914
$ws.Range("D18").Value = "Complete"
$ws.Range("G18").Value = 45914

# Row 19: In Progress -> Complete
$ws.Range("D19").Value = "Complete"

# Row 20: In Progress -> Complete, start/end move
$ws.Range("D20").Value = "Complete"
$ws.Range("F20").Value = 45914
$ws.Range("G20").Value = 45915

# Row 21: In Progress -> Complete
$ws.Range("D21").Value = "Complete"

# Row 22: In Progress -> Complete, start/end move
$ws.Range("D22").Value = "Complete"
$ws.Range("F22").Value = 45914
$ws.Range("G22").Value = 45916

# Row 23: In Progress -> Complete, start/end move
$ws.Range("D23").Value = "Complete"
$ws.Range("F23").Value = 45914
$ws.Range("G23").Value = 45917

# Row 25: Not Started -> Complete, end date moves
$ws.Range("D25").Value = "Complete"
$ws.Range("G25").Value = 45919

# Row 26: Not Started -> Complete, end date moves
$ws.Range("D26").Value = "Complete"
$ws.Range("G26").Value = 45919

# Row 27: Not Started -> In Progress, end date moves
$ws.Range("D27").Value = "In Progress"
$ws.Range("G27").Value = 45929

# Row 28: Not Started -> In Progress, end date moves
$ws.Range("D28").Value = "In Progress"
$ws.Range("G28").Value = 45929

# Row 29: Not Started -> In Progress
$ws.Range("D29").Value = "In Progress"

# Row 30: Not Started -> In Progress
$ws.Range("D30").Value = "In Progress"

# Row 31: Not Started -> In Progress
$ws.Range("D31").Value = "In Progress"

# Row 32: Not Started -> In Progress
$ws.Range("D32").Value = "In Progress"

# Row 33: start date slips
$ws.Range("F33").Value = 45930

# Row 34: start date slips
$ws.Range("F34").Value = 45933

# Row 35: start date slips
$ws.Range("F35").Value = 45933

# ---------------------------------------------------------------------------
# Column E is now hidden
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).Hidden = $true

# ---------------------------------------------------------------------------
# Sheet view: scroll position, zoom and selection changed
# ---------------------------------------------------------------------------
$ws.Activate()
$window = $excel.ActiveWindow
$window.FreezePanes = $false
$window.ScrollColumn = 1
$window.ScrollRow = 1
$ws.Range("A2").Select()
$window.FreezePanes = $true
$window.Zoom = 100
$ws.Range("C37").Select()
$window.ScrollColumn = 3
